$wb = $excel.ActiveWorkbook

# Update "想去人数" (interested count) values for two events that appear
# on both the "展览" sheet and the "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F7").Value = 6266
    $ws.Range("F15").Value = 465
}
